$wb = $excel.ActiveWorkbook

# --- Step 1: clear header-row formatting (bold, border, center/top align) on all sheets ---
foreach ($name in @("ALC","ARM","BSM","CRP","CUL","GSM","LTW","WVR")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A1:N1").ClearFormats()
}

# --- Step 2: apply updated currentAveragePrice / Leve profit figures (refreshed market data) ---

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 36375.74
$ws.Range("I106").Value = 2409.2856
$ws.Range("K106").Value = 2409.2856
$ws.Range("M106").Value = -1778.2856
$ws.Range("H132").Value = 23451.705
$ws.Range("I132").Value = 3198.861
$ws.Range("J132").Value = 114589.5
$ws.Range("K132").Value = 9596.582999999999
$ws.Range("L132").Value = 343768.5
$ws.Range("M132").Value = -7066.582999999999
$ws.Range("N132").Value = -348828.5
$ws.Range("H137").Value = 6244.2793
$ws.Range("I137").Value = 9114.733
$ws.Range("J137").Value = 4706.5356
$ws.Range("K137").Value = 27344.199
$ws.Range("L137").Value = 14119.6068
$ws.Range("M137").Value = -24794.199
$ws.Range("N137").Value = -19219.6068
$ws.Range("H138").Value = 2525.566
$ws.Range("I138").Value = 1460.8125
$ws.Range("J138").Value = 2986
$ws.Range("K138").Value = 4382.4375
$ws.Range("L138").Value = 8958
$ws.Range("M138").Value = 757.5625
$ws.Range("N138").Value = -19238

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11498.186
$ws.Range("I32").Value = 10776.911
$ws.Range("J32").Value = 15104.556
$ws.Range("K32").Value = 10776.911
$ws.Range("L32").Value = 15104.556
$ws.Range("M32").Value = -10489.911
$ws.Range("N32").Value = -15678.556

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1214.9445
$ws.Range("I94").Value = 1241.3572
$ws.Range("K94").Value = 1241.3572
$ws.Range("M94").Value = -790.3571999999999
$ws.Range("H99").Value = 1916.9166
$ws.Range("I99").Value = 1800.6552
$ws.Range("K99").Value = 1800.6552
$ws.Range("M99").Value = -302.6551999999999
$ws.Range("H105").Value = 3306.65
$ws.Range("I105").Value = 1435.7142
$ws.Range("J105").Value = 4314.077
$ws.Range("K105").Value = 1435.7142
$ws.Range("L105").Value = 4314.077
$ws.Range("M105").Value = 311.2858000000001
$ws.Range("N105").Value = -7808.077

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 31379538
$ws.Range("I4").Value = 1000000000
$ws.Range("J4").Value = 133717.03
$ws.Range("K4").Value = 1000000000
$ws.Range("L4").Value = 133717.03
$ws.Range("M4").Value = -999999888
$ws.Range("N4").Value = -133941.03
$ws.Range("H10").Value = 174.5
$ws.Range("I10").Value = 174.5
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 174.5
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -35.5
$ws.Range("N10").ClearContents()
$ws.Range("H21").Value = 16000
$ws.Range("J21").Value = 16000
$ws.Range("L21").Value = 16000
$ws.Range("N21").Value = -16470
$ws.Range("H22").Value = 1247.6666
$ws.Range("I22").Value = 294.9
$ws.Range("J22").Value = 2113.818
$ws.Range("K22").Value = 294.9
$ws.Range("L22").Value = 2113.818
$ws.Range("M22").Value = 55.10000000000002
$ws.Range("N22").Value = -2813.818
$ws.Range("H23").Value = 166704850
$ws.Range("I23").Value = 333339680
$ws.Range("K23").Value = 333339680
$ws.Range("M23").Value = -333339440
$ws.Range("H26").Value = 22394.955
$ws.Range("J26").Value = 22394.955
$ws.Range("L26").Value = 22394.955
$ws.Range("N26").Value = -22968.955
$ws.Range("H27").Value = 166704850
$ws.Range("I27").Value = 333339680
$ws.Range("K27").Value = 333339680
$ws.Range("M27").Value = -333339488
$ws.Range("H32").Value = 1603.6666
$ws.Range("I32").Value = 905.5
$ws.Range("J32").Value = 3000
$ws.Range("K32").Value = 905.5
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -589.5
$ws.Range("N32").Value = -3632
$ws.Range("H36").Value = 26335.572
$ws.Range("I36").Value = 1299.5
$ws.Range("J36").Value = 36350
$ws.Range("K36").Value = 1299.5
$ws.Range("L36").Value = 36350
$ws.Range("M36").Value = -911.5
$ws.Range("N36").Value = -37126
$ws.Range("H39").Value = 12499.75
$ws.Range("I39").Value = 2666.6667
$ws.Range("J39").Value = 41999
$ws.Range("K39").Value = 2666.6667
$ws.Range("L39").Value = 41999
$ws.Range("M39").Value = -2275.6667
$ws.Range("N39").Value = -42781
$ws.Range("H40").Value = 26335.572
$ws.Range("I40").Value = 1299.5
$ws.Range("J40").Value = 36350
$ws.Range("K40").Value = 1299.5
$ws.Range("L40").Value = 36350
$ws.Range("M40").Value = -1139.5
$ws.Range("N40").Value = -36670
$ws.Range("H44").Value = 34580.145
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H49").Value = 12499.75
$ws.Range("I49").Value = 2666.6667
$ws.Range("J49").Value = 41999
$ws.Range("K49").Value = 2666.6667
$ws.Range("L49").Value = 41999
$ws.Range("M49").Value = -2484.6667
$ws.Range("N49").Value = -42363
$ws.Range("H58").Value = 10205371
$ws.Range("I58").Value = 774.7059
$ws.Range("J58").Value = 33335788
$ws.Range("K58").Value = 774.7059
$ws.Range("L58").Value = 33335788
$ws.Range("M58").Value = -571.7059
$ws.Range("N58").Value = -33336194
$ws.Range("H86").Value = 4191.2
$ws.Range("J86").Value = 2441
$ws.Range("L86").Value = 2441
$ws.Range("N86").Value = -4687
$ws.Range("H89").Value = 4191.2
$ws.Range("J89").Value = 2441
$ws.Range("L89").Value = 12205
$ws.Range("N89").Value = -23437
$ws.Range("H136").Value = 10205371
$ws.Range("I136").Value = 774.7059
$ws.Range("J136").Value = 33335788
$ws.Range("K136").Value = 2324.1177
$ws.Range("L136").Value = 100007364
$ws.Range("M136").Value = 225.8822999999998
$ws.Range("N136").Value = -100012464

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5387.5
$ws.Range("I113").Value = 8422.538
$ws.Range("J113").Value = 1003.55554
$ws.Range("K113").Value = 25267.614
$ws.Range("L113").Value = 3010.66662
$ws.Range("M113").Value = -23097.614
$ws.Range("N113").Value = -7350.66662

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 20809.096
$ws.Range("I126").Value = 37522.453
$ws.Range("J126").Value = 2424.4
$ws.Range("K126").Value = 112567.359
$ws.Range("L126").Value = 7273.200000000001
$ws.Range("M126").Value = -110097.359
$ws.Range("N126").Value = -12213.2
$ws.Range("H132").Value = 2543.7437
$ws.Range("I132").Value = 1826.9131
$ws.Range("J132").Value = 3574.1875
$ws.Range("K132").Value = 5480.7393
$ws.Range("L132").Value = 10722.5625
$ws.Range("M132").Value = -2950.7393
$ws.Range("N132").Value = -15782.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1170.0588
$ws.Range("I93").Value = 1048
$ws.Range("K93").Value = 1048
$ws.Range("M93").Value = 200
